# Fruta / hortaliza, semanal
# Insert a new weekly record as row 27, shifting the existing rows 27-34
# down to rows 28-35 (matching the target dimension A1:R35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 27.
$ws.Rows.Item(27).Insert()

# Populate the new row 27 with the new weekly market record.
$ws.Cells.Item(27, 1).Value2  = 3
$ws.Cells.Item(27, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(27, 3).Value2  = "Coquimbo"
$ws.Cells.Item(27, 4).Value2  = 44455
$ws.Cells.Item(27, 5).Value2  = 5
$ws.Cells.Item(27, 6).Value2  = 100112035
$ws.Cells.Item(27, 7).Value2  = "Bruselas (repollito)"
$ws.Cells.Item(27, 8).Value2  = "Sin especificar"
$ws.Cells.Item(27, 9).Value2  = "Primera"
$ws.Cells.Item(27, 10).Value2 = 35
$ws.Cells.Item(27, 11).Value2 = 22000
$ws.Cells.Item(27, 12).Value2 = 22000
$ws.Cells.Item(27, 13).Value2 = 22000
$ws.Cells.Item(27, 14).Value2 = "$/malla 15 kilos"
$ws.Cells.Item(27, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(27, 16).Value2 = 1467
$ws.Cells.Item(27, 17).Value2 = 15
$ws.Cells.Item(27, 18).Value2 = "Hortaliza"
